# Scheduled-runner market data refresh for Brynhildr_Profits.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for the
# leves whose underlying item prices moved, across all 8 crafting-job sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Values are written as literals
# (no formulas exist in these tables) to mirror the upstream price-bot output.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1237
$ws.Range("I19").Value = 1147.75
$ws.Range("J19").Value = 1326.25
$ws.Range("K19").Value = 1147.75
$ws.Range("L19").Value = 1326.25
$ws.Range("M19").Value = -972.75
$ws.Range("N19").Value = -1676.25
# Row 33
$ws.Range("H33").Value = 210.86667
$ws.Range("I33").Value = 172.08333
$ws.Range("K33").Value = 172.08333
$ws.Range("M33").Value = 56.91667000000001
# Row 38
$ws.Range("H38").Value = 1449
$ws.Range("I38").Value = 612.5454999999999
$ws.Range("K38").Value = 1837.6365
$ws.Range("M38").Value = -1465.6365
# Row 61
$ws.Range("H61").Value = 2033.3334
$ws.Range("I61").Value = 1000
$ws.Range("J61").Value = 2550
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 7650
$ws.Range("M61").Value = -2828
$ws.Range("N61").Value = -7994
# Row 70
$ws.Range("H70").Value = 2826.889
$ws.Range("I70").Value = 2162.75
$ws.Range("J70").Value = 3358.2
$ws.Range("K70").Value = 6488.25
$ws.Range("L70").Value = 10074.6
$ws.Range("M70").Value = -6218.25
$ws.Range("N70").Value = -10614.6
# Row 73
$ws.Range("H73").Value = 2826.889
$ws.Range("I73").Value = 2162.75
$ws.Range("J73").Value = 3358.2
$ws.Range("K73").Value = 6488.25
$ws.Range("L73").Value = 10074.6
$ws.Range("M73").Value = -5552.25
$ws.Range("N73").Value = -11946.6
# Row 86
$ws.Range("H86").Value = 3900
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 2800
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 2800
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -5046
# Row 89
$ws.Range("H89").Value = 3900
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 2800
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 14000
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -25232
# Row 103
$ws.Range("H103").Value = 623.78723
$ws.Range("I103").Value = 508.95
$ws.Range("K103").Value = 1526.85
$ws.Range("M103").Value = -940.8499999999999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 168630.48
$ws.Range("I32").Value = 173915.92
$ws.Range("K32").Value = 173915.92
$ws.Range("M32").Value = -173628.92
# Row 74
$ws.Range("H74").Value = 4815.8823
$ws.Range("I74").Value = 739.2368
$ws.Range("K74").Value = 739.2368
$ws.Range("M74").Value = 134.7632
# Row 77
$ws.Range("H77").Value = 4815.8823
$ws.Range("I77").Value = 739.2368
$ws.Range("K77").Value = 3696.184
$ws.Range("M77").Value = 671.8159999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 29
$ws.Range("H29").Value = 500
$ws.Range("I29").Value = 500
$ws.Range("K29").Value = 500
$ws.Range("M29").Value = -211
# Row 99
$ws.Range("H99").Value = 6761.4443
$ws.Range("I99").Value = 7431.625
$ws.Range("K99").Value = 7431.625
$ws.Range("M99").Value = -5933.625
# Row 107
$ws.Range("H107").Value = 1014.8571
$ws.Range("I107").Value = 934
$ws.Range("K107").Value = 934
$ws.Range("M107").Value = 986

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 91776.27
$ws.Range("I16").Value = 886.125
$ws.Range("K16").Value = 886.125
$ws.Range("M16").Value = -599.125
# Row 31
$ws.Range("H31").Value = 2249.7256
$ws.Range("I31").Value = 1901.2632
$ws.Range("K31").Value = 1901.2632
$ws.Range("M31").Value = -1606.2632
# Row 34
$ws.Range("H34").Value = 2249.7256
$ws.Range("I34").Value = 1901.2632
$ws.Range("K34").Value = 1901.2632
$ws.Range("M34").Value = -1699.2632
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
# Row 86
$ws.Range("H86").Value = 17909.818
$ws.Range("I86").Value = 19889.777
$ws.Range("K86").Value = 19889.777
$ws.Range("M86").Value = -18766.777
# Row 89
$ws.Range("H89").Value = 17909.818
$ws.Range("I89").Value = 19889.777
$ws.Range("K89").Value = 99448.88499999999
$ws.Range("M89").Value = -93832.88499999999
# Row 105
$ws.Range("H105").Value = 8279.143
$ws.Range("I105").Value = 10601
$ws.Range("J105").Value = 2474.5
$ws.Range("K105").Value = 10601
$ws.Range("L105").Value = 2474.5
$ws.Range("M105").Value = -8854
$ws.Range("N105").Value = -5968.5
# Row 113
$ws.Range("H113").Value = 91776.27
$ws.Range("I113").Value = 886.125
$ws.Range("K113").Value = 886.125
$ws.Range("M113").Value = 1283.875
# Row 134
$ws.Range("H134").Value = 1415.4
$ws.Range("I134").Value = 1419.8889
$ws.Range("K134").Value = 4259.6667
$ws.Range("M134").Value = -1724.6667

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 12727416
$ws.Range("I4").Value = 15555642
$ws.Range("K4").Value = 46666926
$ws.Range("M4").Value = -46666814
# Row 68
$ws.Range("H68").Value = 6363.364
$ws.Range("J68").Value = 6363.364
$ws.Range("L68").Value = 19090.092
$ws.Range("N68").Value = -20712.092
# Row 71
$ws.Range("H71").Value = 6363.364
$ws.Range("J71").Value = 6363.364
$ws.Range("L71").Value = 57270.276
$ws.Range("N71").Value = -65382.276
# Row 75
$ws.Range("H75").Value = 2068.6667
$ws.Range("J75").Value = 2606
$ws.Range("L75").Value = 7818
$ws.Range("N75").Value = -9814
# Row 78
$ws.Range("H78").Value = 2068.6667
$ws.Range("J78").Value = 2606
$ws.Range("L78").Value = 23454
$ws.Range("N78").Value = -33438
# Row 92
$ws.Range("H92").Value = 745
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
# Row 94
$ws.Range("H94").Value = 5997.5
$ws.Range("I94").Value = 5997.5
$ws.Range("K94").Value = 17992.5
$ws.Range("M94").Value = -17316.5
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
# Row 129
$ws.Range("H129").Value = 1003336.2
$ws.Range("J129").Value = 3104.7144
$ws.Range("L129").Value = 9314.143199999999
$ws.Range("N129").Value = -19314.1432
# Row 133
$ws.Range("H133").Value = 5462
$ws.Range("I133").Value = 5405.3687
$ws.Range("J133").Value = 6000
$ws.Range("K133").Value = 16216.1061
$ws.Range("L133").Value = 18000
$ws.Range("M133").Value = -11156.1061
$ws.Range("N133").Value = -28120
# Row 137
$ws.Range("H137").Value = 8181.5
$ws.Range("I137").Value = 5030
$ws.Range("J137").Value = 11333
$ws.Range("K137").Value = 15090
$ws.Range("L137").Value = 33999
$ws.Range("M137").Value = -9990
$ws.Range("N137").Value = -44199
# Row 139
$ws.Range("H139").Value = 4017.1428
$ws.Range("I139").Value = 2572.2
$ws.Range("K139").Value = 7716.599999999999
$ws.Range("M139").Value = -2576.599999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 8793.014999999999
$ws.Range("I132").Value = 12084.381
$ws.Range("J132").Value = 3673.111
$ws.Range("K132").Value = 36253.143
$ws.Range("L132").Value = 11019.333
$ws.Range("M132").Value = -33723.143
$ws.Range("N132").Value = -16079.333
# Row 134
$ws.Range("H134").Value = 51062.5
$ws.Range("J134").Value = 51062.5
$ws.Range("L134").Value = 153187.5
$ws.Range("N134").Value = -158257.5
# Row 136
$ws.Range("H136").Value = 73665.164
$ws.Range("J136").Value = 73665.164
$ws.Range("L136").Value = 220995.492
$ws.Range("N136").Value = -226095.492

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3018
$ws.Range("I40").Value = 1822
$ws.Range("J40").Value = 4413.3335
$ws.Range("K40").Value = 1822
$ws.Range("L40").Value = 4413.3335
$ws.Range("M40").Value = -1686
$ws.Range("N40").Value = -4685.3335

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 24
$ws.Range("H24").Value = 349666.66
$ws.Range("I24").Value = 502000
$ws.Range("J24").Value = 45000
$ws.Range("K24").Value = 502000
$ws.Range("L24").Value = 45000
$ws.Range("M24").Value = -501770
$ws.Range("N24").Value = -45460
# Row 29
$ws.Range("H29").Value = 100000
$ws.Range("I29").Value = 100000
$ws.Range("K29").Value = 100000
$ws.Range("M29").Value = -99710
# Row 62
$ws.Range("H62").Value = 4999.75
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
# Row 65
$ws.Range("H65").Value = 4999.75
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380
# Row 107
$ws.Range("H107").Value = 1584.8823
$ws.Range("J107").Value = 3910
$ws.Range("L107").Value = 11730
$ws.Range("N107").Value = -15570
# Row 121
$ws.Range("H121").Value = 60000
$ws.Range("J121").Value = 60000
$ws.Range("L121").Value = 60000
$ws.Range("N121").Value = -63494
# Row 122
$ws.Range("H122").Value = 33186.027
$ws.Range("I122").Value = 1849.625
$ws.Range("J122").Value = 95858.836
$ws.Range("K122").Value = 5548.875
$ws.Range("L122").Value = 287576.508
$ws.Range("M122").Value = -3098.875
$ws.Range("N122").Value = -292476.508
# Row 123
$ws.Range("H123").Value = 59999.168
$ws.Range("J123").Value = 59999.168
$ws.Range("L123").Value = 59999.168
$ws.Range("N123").Value = -69799.16800000001
# Row 136
$ws.Range("H136").Value = 1057.0646
$ws.Range("I136").Value = 1166.6522
$ws.Range("J136").Value = 742
$ws.Range("K136").Value = 3499.9566
$ws.Range("L136").Value = 2226
$ws.Range("M136").Value = -949.9566
$ws.Range("N136").Value = -7326
